$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Meter place" header column (S) and its values.
$ws.Range("S1").Value = "Meter place"

$ws.Range("S2").Value = "Kitchen"
$ws.Range("S3").Value = "Bathroom"
$ws.Range("S4").Value = "Kitchen"
$ws.Range("S5").Value = "Bathroom"
$ws.Range("S6").Value = "Kitchen"
$ws.Range("S7").Value = "Bathroom"
$ws.Range("S8").Value = "Kitchen"
$ws.Range("S9").Value = "Bathroom"
$ws.Range("S10").Value = "Kitchen"
$ws.Range("S11").Value = "Bathroom"

# Match column S width/style to column R.
$ws.Range("S1:S11").Style = $ws.Range("R1:R11").Style
$ws.Columns.Item(19).ColumnWidth = $ws.Columns.Item(18).ColumnWidth
